# Update the "想去人数" (F column) counters on the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1269
$ws1.Range("F3").Value  = 664
$ws1.Range("F4").Value  = 361
$ws1.Range("F5").Value  = 5122
$ws1.Range("F6").Value  = 547
$ws1.Range("F7").Value  = 9892
$ws1.Range("F8").Value  = 255
$ws1.Range("F9").Value  = 548
$ws1.Range("F10").Value = 97
$ws1.Range("F11").Value = 43
$ws1.Range("F12").Value = 741

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1269
$ws4.Range("F3").Value  = 664
$ws4.Range("F4").Value  = 361
$ws4.Range("F7").Value  = 5122
$ws4.Range("F8").Value  = 547
$ws4.Range("F10").Value = 9892
$ws4.Range("F11").Value = 255
$ws4.Range("F12").Value = 548
$ws4.Range("F13").Value = 97
$ws4.Range("F16").Value = 43
$ws4.Range("F17").Value = 741
$ws4.Range("F18").Value = 1
